# Added games for 1/20/2021
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in "Beat Vegas?" results for the existing 1/19/2021 games (rows 97-98) ---
$ws.Range("G97").Value = "No"
$ws.Range("G98").Value = "No"

# --- New games for 1/20/2021 (rows 99-107) ---
$newGames = @(
    @{ Row = 99;  Date = 44216; Home = "CLE"; Away = "BRK"; Spread = 10;    Pred = 30.4;               Diff = -20.399999999999999 },
    @{ Row = 100; Date = 44216; Home = "PHI"; Away = "BOS"; Spread = -5;    Pred = -5.6;               Diff = 0.59999999999999964 },
    @{ Row = 101; Date = 44216; Home = "IND"; Away = "DAL"; Spread = 1.5;  Pred = -6.7;               Diff = 8.1999999999999993 },
    @{ Row = 102; Date = 44216; Home = "ATL"; Away = "DET"; Spread = -5;    Pred = -2.5;               Diff = -2.5 },
    @{ Row = 103; Date = 44216; Home = "TOR"; Away = "MIA"; Spread = -4.5;  Pred = 8.8000000000000007; Diff = -13.3 },
    @{ Row = 104; Date = 44216; Home = "MIN"; Away = "ORL"; Spread = 4;    Pred = -5.3;               Diff = 9.3000000000000007 },
    @{ Row = 105; Date = 44216; Home = "HOU"; Away = "PHO"; Spread = 5.5;  Pred = -0.2;               Diff = 5.7 },
    @{ Row = 106; Date = 44216; Home = "LAC"; Away = "SAC"; Spread = -9.5;  Pred = -4;                 Diff = -5.5 },
    @{ Row = 107; Date = 44216; Home = "GSW"; Away = "SAS"; Spread = -1;    Pred = -4.2;               Diff = 3.2 }
)

foreach ($g in $newGames) {
    $r = $g.Row
    $ws.Range("A$r").Value = $g.Date
    $ws.Range("A$r").NumberFormat = "yyyy\-mm\-dd"
    $ws.Range("B$r").Value = $g.Home
    $ws.Range("C$r").Value = $g.Away
    $ws.Range("D$r").Value = $g.Spread
    $ws.Range("E$r").Value = $g.Pred
    $ws.Range("F$r").Value = $g.Diff
}

# --- Update view state to match where the author left off editing ---
$ws.Range("J107").Select()
